$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add two new bullet items under "Errores:" after the paragraph
#    "Cuando el admin quiere obtener las salas tira un error el back":
#       - "La flor no está sumando puntos"
#       - "Al terminar la partida" + " no está sumando puntos" (two runs)
# ---------------------------------------------------------------------------

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Cuando el admin quiere obtener las salas*") {
        $anchor = $p.Range
        break
    }
}

$pos = $anchor.End
# New empty paragraph, inherits the same pStyle/numPr as the anchor paragraph.
$anchor.InsertParagraphAfter() | Out-Null

$t1 = "La flor no está sumando puntos"
$d.Range($pos, $pos).InsertBefore($t1) | Out-Null
$pos = $pos + $t1.Length

# Second new (still empty) list paragraph.
$d.Range($pos, $pos).InsertParagraphAfter() | Out-Null
$pos = $pos + 1

# Fill it with two distinct runs via raw OOXML so they don't get merged into one.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$p2Range = $d.Range($pos, $pos + 1)
$p2Xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Prrafodelista"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t>Al terminar la partida</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> no está sumando puntos</w:t></w:r>' +
  '</w:p>' +
  '</pkg:xmlData>'
$p2Range.InsertXML($p2Xml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the second
#    "Truco + Retruco + VC" row to the second "Truco + Retruco" row.
# ---------------------------------------------------------------------------

$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $label = $tbl.Cell($r, 1).Range.Text
    $response = $tbl.Cell($r, 2).Range.Text

    if ($label -like "Truco + Retruco`r*" -and $response -like "No`r*") {
        $cellRange = $tbl.Cell($r, 1).Range
        $xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<w:p ' + $wNs + '><w:r><w:lastRenderedPageBreak/><w:t>Truco + Retruco</w:t></w:r></w:p>' +
          '</pkg:xmlData>'
        $cellRange.InsertXML($xml) | Out-Null
    }
    elseif ($label -like "Truco + Retruco + VC`r*" -and $response -like "No`r*") {
        $cellRange = $tbl.Cell($r, 1).Range
        $xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<w:p ' + $wNs + '><w:r><w:t>Truco + Retruco + VC</w:t></w:r></w:p>' +
          '</pkg:xmlData>'
        $cellRange.InsertXML($xml) | Out-Null
    }
}

Write-Output "edit complete"
